$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (B6:B9) that will be replaced
$ws.Range("B6:B9").ClearContents()

# New headers (order matters for shared string table insertion order)
$ws.Range("E2").Value = "Min"
$ws.Range("G2").Value = "Count"
$ws.Range("F2").Value = "Max"

# New data values B4:B7
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.2
$ws.Range("B6").Value = 0.3
$ws.Range("B7").Value = 0.4

# Formulas row 3 (summary row)
$ws.Range("C3").Formula = "=_xlfn.STDEV.P(B:B)"
$ws.Range("D3").Formula = "=SUM(B:B)"
$ws.Range("E3").Formula = "=MIN(B:B)"
$ws.Range("F3").Formula = "=MAX(B:B)"
$ws.Range("G3").Formula = "=COUNT(B:B)"

# Formulas row 4 (duplicate of row 3 results)
$ws.Range("C4").Formula = "=C3"
$ws.Range("D4").Formula = "=D3"
$ws.Range("E4").Formula = "=E3"
$ws.Range("F4").Formula = "=F3"

# Apply header style (bold) to new header cells
$ws.Range("E2").Font.Bold = $true
$ws.Range("F2").Font.Bold = $true
$ws.Range("G2").Font.Bold = $true

# Apply number formats matching existing columns
$ws.Range("E3").NumberFormat = "0.00000000000000000E+00"
$ws.Range("F3").NumberFormat = "0.00000000000000000E+00"
$ws.Range("E4").NumberFormat = "0.00000000000000000"
$ws.Range("F4").NumberFormat = "0.00000000000000000"

# Column widths
$ws.Range("C1:F1").EntireColumn.ColumnWidth = 33.83
$ws.Range("G1").EntireColumn.ColumnWidth = 19.33

# Selection
$ws.Range("B4").Select()
